# Applies the "Output/Animal10.xlsx" protocol-change edits:
#  - rename header E1 from "strength (raw)" to "strength (RMS)"
#  - update reactionTime (B), difference (D) and strength (E) columns
#    for data rows 2..19 to reflect the new RMS-based computation
#    (peakTime column C is unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename
$ws.Range("E1").Value = "strength (RMS)"

# New values per row: reactionTime, difference, strength
$data = @{
    2  = @{ B = 14;    D = 13.6;  E = 76.2 }
    3  = @{ B = 14.5;  D = 16.5;  E = 79.5 }
    4  = @{ B = 14;    D = 16.8;  E = 75.2 }
    5  = @{ B = 14;    D = 21.2;  E = 77.8 }
    6  = @{ B = 14;    D = 16.4;  E = 73.2 }
    7  = @{ B = 14.4;  D = 14;    E = 67.2 }
    8  = @{ B = 14;    D = 16.4;  E = 72.6 }
    9  = @{ B = 14;    D = 18.8;  E = 73 }
    10 = @{ B = 12.8;  D = 19.6;  E = 77.4 }
    11 = @{ B = 14;    D = 31.6;  E = 71.8 }
    12 = @{ B = 14.4;  D = 16.8;  E = 66.8 }
    13 = @{ B = 14;    D = 17.6;  E = 71.4 }
    14 = @{ B = 13.33; D = 17.33; E = 73.67 }
    15 = @{ B = 14;    D = 17;    E = 66.25 }
    16 = @{ B = 14;    D = 20.4;  E = 78.8 }
    17 = @{ B = 13.6;  D = 20.8;  E = 82.2 }
    18 = @{ B = 15.2;  D = 16;    E = 69.6 }
    19 = @{ B = 12.8;  D = 18;    E = 72.8 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 5).Value = $vals.E
}
